$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:E2").NumberFormat = "@"
$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '30.124.35'
$ws.Range("E2").Value = '  +0.95%  '

$ws.Range("B3:E3").NumberFormat = "@"
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.891.88'
$ws.Range("E3").Value = '  +0.19%  '

$ws.Range("B4:E4").NumberFormat = "@"
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("B5:E5").NumberFormat = "@"
$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").Value = '0.7393'
$ws.Range("E5").Value = '  -1.38%  '

$ws.Range("B6:E6").NumberFormat = "@"
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '242.64'
$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("B7:E7").NumberFormat = "@"
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("B8:E8").NumberFormat = "@"
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.3173'
$ws.Range("E8").Value = '  +1.42%  '

$ws.Range("B9:E9").NumberFormat = "@"
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.07207'
$ws.Range("E9").Value = '  +1.20%  '

$ws.Range("B10:E10").NumberFormat = "@"
$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").Value = '24.93'
$ws.Range("E10").Value = '  -1.36%  '

$ws.Range("B11:E11").NumberFormat = "@"
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '0.08345'
$ws.Range("E11").Value = '  -1.76%  '

$ws.Range("B12:E12").NumberFormat = "@"
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = '0.7605'
$ws.Range("E12").Value = '  +0.03%  '

$ws.Range("B13:E13").NumberFormat = "@"
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.924.52'
$ws.Range("E13").Value = '  +1.08%  '

$ws.Range("B14:E14").NumberFormat = "@"
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.451'
$ws.Range("E14").Value = '  +1.64%  '

$ws.Range("B15:E15").NumberFormat = "@"
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '93.15'
$ws.Range("E15").Value = '  -0.27%  '

$ws.Range("B16:E16").NumberFormat = "@"
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '30.254.35'
$ws.Range("E16").Value = '  +2.04%  '

$ws.Range("B17:E17").NumberFormat = "@"
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = '6.175'
$ws.Range("E17").Value = '  +0.39%  '

$ws.Range("B18:E18").NumberFormat = "@"
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").Value = '250.59'
$ws.Range("E18").Value = '  +2.88%  '

$ws.Range("B19:E19").NumberFormat = "@"
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = '13.64'
$ws.Range("E19").Value = '  -0.50%  '

$ws.Range("B20:E20").NumberFormat = "@"
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.000007895'
$ws.Range("E20").Value = '  +1.19%  '

$ws.Range("B21:E21").NumberFormat = "@"
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.179.17'
$ws.Range("E21").Value = '  +1.21%  '

$ws.Range("B22:E22").NumberFormat = "@"
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '0.9993'
$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("B23:E23").NumberFormat = "@"
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = '7.954'
$ws.Range("E23").Value = '  -0.25%  '

$ws.Range("B24:E24").NumberFormat = "@"
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").Value = '0.9995'
$ws.Range("E24").Value = '  -0.10%  '

$ws.Range("B25:E25").NumberFormat = "@"
$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D25").Value = '0.1584'
$ws.Range("E25").Value = '  -0.40%  '

$ws.Range("B26:E26").NumberFormat = "@"
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '9.318'
$ws.Range("E26").Value = '  -0.57%  '

$ws.Range("B27:E27").NumberFormat = "@"
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '164.64'
$ws.Range("E27").Value = '  +1.36%  '

$ws.Range("B28:E28").NumberFormat = "@"
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '18.79'
$ws.Range("E28").Value = '  +0.28%  '

$ws.Range("B29:E29").NumberFormat = "@"
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '2.066'
$ws.Range("E29").Value = '  +1.83%  '

$ws.Range("B30:E30").NumberFormat = "@"
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '1.488'
$ws.Range("E30").Value = '  -0.41%  '

$ws.Range("B31:E31").NumberFormat = "@"
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '4.589'
$ws.Range("E31").Value = '  +1.47%  '

$ws.Range("B32:E32").NumberFormat = "@"
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '1.539'
$ws.Range("E32").Value = '  +0.09%  '

$ws.Range("B33:E33").NumberFormat = "@"
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '4.208'
$ws.Range("E33").Value = '  +1.92%  '

$ws.Range("B34:E34").NumberFormat = "@"
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.05372'
$ws.Range("E34").Value = '  -0.93%  '

$ws.Range("B35:E35").NumberFormat = "@"
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.258'
$ws.Range("E35").Value = '  +1.43%  '

$ws.Range("B36:E36").NumberFormat = "@"
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.7778'
$ws.Range("E36").Value = '  +4.14%  '

$ws.Range("B37:E37").NumberFormat = "@"
$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").Value = '1.004'
$ws.Range("E37").Value = '  -0.13%  '

$ws.Range("B38:E38").NumberFormat = "@"
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").Value = '2.728'
$ws.Range("E38").Value = '  +0.65%  '

$ws.Range("B39:E39").NumberFormat = "@"
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01964'
$ws.Range("E39").Value = '  +0.85%  '

$ws.Range("B40:E40").NumberFormat = "@"
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.765'
$ws.Range("E40").Value = '  -0.33%  '

$ws.Range("B41:E41").NumberFormat = "@"
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.4570'
$ws.Range("E41").Value = '  +2.56%  '

$ws.Range("B42:E42").NumberFormat = "@"
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '6.103'
$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("B43:E43").NumberFormat = "@"
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.100.36'
$ws.Range("E43").Value = '  +0.47%  '

$ws.Range("B44:E44").NumberFormat = "@"
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '72.69'
$ws.Range("E44").Value = '  +0.21%  '

$ws.Range("B45:E45").NumberFormat = "@"
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '0.8774'
$ws.Range("E45").Value = '  +2.43%  '

$ws.Range("B46:E46").NumberFormat = "@"
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '104.67'
$ws.Range("E46").Value = '  +2.44%  '

$ws.Range("B47:E47").NumberFormat = "@"
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.10%  '

$ws.Range("B48:E48").NumberFormat = "@"
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '1.867'
$ws.Range("E48").Value = '  +0.28%  '

$ws.Range("B49:E49").NumberFormat = "@"
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").Value = '7.602'
$ws.Range("E49").Value = '  -1.36%  '

$ws.Range("B50:E50").NumberFormat = "@"
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.082.09'
$ws.Range("E50").Value = '  +1.80%  '

$ws.Range("B51:E51").NumberFormat = "@"
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '9.632'
$ws.Range("E51").Value = '  -1.19%  '
